$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3419041688096627
$ws.Range("D2").Value = 0.3791904208309234
$ws.Range("E2").Value = 0.1143888469716536
$ws.Range("F2").Value = 5.241275693317704
$ws.Range("G2").Value = 0.002644826841643954
$ws.Range("K2").Value = 2.736518024821009
$ws.Range("L2").Value = 0.07823533792780779
$ws.Range("M2").Value = 0.7778433907156881

$ws.Range("C3").Value = 0.3384035684953233
$ws.Range("D3").Value = 0.3698314667457367
$ws.Range("E3").Value = 0.1136827103935509
$ws.Range("F3").Value = 5.059717607612498
$ws.Range("G3").Value = 0.002652081424053018
$ws.Range("K3").Value = 2.658006374373485
$ws.Range("L3").Value = 0.07799861627761118
$ws.Range("M3").Value = 0.7618770079938457

$ws.Range("C4").Value = 0.3364207824093484
$ws.Range("D4").Value = 0.3641025473780672
$ws.Range("E4").Value = 0.1132923073987691
$ws.Range("F4").Value = 4.949714129688232
$ws.Range("G4").Value = 0.002656761340918131
$ws.Range("K4").Value = 2.612415698208508
$ws.Range("L4").Value = 0.07787153359670995
$ws.Range("M4").Value = 0.7527096143451359

$ws.Range("C5").Value = 0.335654499375309
$ws.Range("D5").Value = 0.3617717749973082
$ws.Range("E5").Value = 0.1131440283278735
$ws.Range("F5").Value = 4.905248316768819
$ws.Range("G5").Value = 0.002658725397887755
$ws.Range("K5").Value = 2.594491320102065
$ws.Range("L5").Value = 0.07782431972290382
$ws.Range("M5").Value = 0.7491330284808058

$ws.Range("C6").Value = 0.3355297733680374
$ws.Range("D6").Value = 0.3613849647175016
$ws.Range("E6").Value = 0.1131200586058192
$ws.Range("F6").Value = 4.897886370903535
$ws.Range("G6").Value = 0.002659054974463435
$ws.Range("K6").Value = 2.591554400173607
$ws.Range("L6").Value = 0.07781675554820211
$ws.Range("M6").Value = 0.7485487335177652

$ws.Range("C7").Value = 0.3364102793570822
$ws.Range("D7").Value = 0.364071099131337
$ws.Range("E7").Value = 0.1132902639264977
$ws.Range("F7").Value = 4.9491129964903
$ws.Range("G7").Value = 0.002656787597994404
$ws.Range("K7").Value = 2.612171319595916
$ws.Range("L7").Value = 0.07787087836211626
$ws.Range("M7").Value = 0.7526607356042447

$ws.Range("C8").Value = 0.3406624769759219
$ws.Range("D8").Value = 0.3759594559080739
$ws.Range("E8").Value = 0.1141363858959821
$ws.Range("F8").Value = 5.178362973687058
$ws.Range("G8").Value = 0.002647281542790264
$ws.Range("K8").Value = 2.708902083185706
$ws.Range("L8").Value = 0.07814991325455978
$ws.Range("M8").Value = 0.7722057179892445

$ws.Range("C9").Value = 0.3503317208017052
$ws.Range("D9").Value = 0.3994393854256231
$ws.Range("E9").Value = 0.1161402446474611
$ws.Range("F9").Value = 5.640054207072581
$ws.Range("G9").Value = 0.002630419482497055
$ws.Range("K9").Value = 2.9195270527415
$ws.Range("L9").Value = 0.07884299498048364
$ws.Range("M9").Value = 0.8156174382566093

$ws.Range("C10").Value = 0.3582606153416066
$ws.Range("D10").Value = 0.4168320012945514
$ws.Range("E10").Value = 0.1178258033005299
$ws.Range("F10").Value = 5.987299170909637
$ws.Range("G10").Value = 0.002619100804996688
$ws.Range("K10").Value = 3.087316499813312
$ws.Range("L10").Value = 0.07944260414475934
$ws.Range("M10").Value = 0.8506696687500437

$ws.Range("C11").Value = 0.3620499260063639
$ws.Range("D11").Value = 0.4247840028010046
$ws.Range("E11").Value = 0.118639668346578
$ws.Range("F11").Value = 6.147158992942082
$ws.Range("G11").Value = 0.002614180772225802
$ws.Range("K11").Value = 3.166544850635375
$ws.Range("L11").Value = 0.07973534205719091
$ws.Range("M11").Value = 0.8673149067272874

$ws.Range("C12").Value = 0.3635113117697415
$ws.Range("D12").Value = 0.4278016707648362
$ws.Range("E12").Value = 0.1189546873177747
$ws.Range("F12").Value = 6.207977602810729
$ws.Range("G12").Value = 0.002612350351099249
$ws.Range("K12").Value = 3.19696859912375
$ws.Range("L12").Value = 0.07984909210198055
$ws.Range("M12").Value = 0.8737196955006965

$ws.Range("C13").Value = 0.3631953959072689
$ws.Range("D13").Value = 0.4271514658954345
$ws.Range("E13").Value = 0.1188865378903863
$ws.Range("F13").Value = 6.19486646445165
$ws.Range("G13").Value = 0.002612743114795235
$ws.Range("K13").Value = 3.190397469730101
$ws.Range("L13").Value = 0.07982446477348759
$ws.Range("M13").Value = 0.8723357752541858

$ws.Range("C14").Value = 0.3621696235620675
$ws.Range("D14").Value = 0.4250321357107225
$ws.Range("E14").Value = 0.1186654480321003
$ws.Range("F14").Value = 6.152156840904183
$ws.Range("G14").Value = 0.002614029528671791
$ws.Range("K14").Value = 3.169039352134234
$ws.Range("L14").Value = 0.07974464214822419
$ws.Range("M14").Value = 0.8678397907081461

$ws.Range("C15").Value = 0.3615447605317286
$ws.Range("D15").Value = 0.4237348403820249
$ws.Range("E15").Value = 0.1185309146673603
$ws.Range("F15").Value = 6.126033153577225
$ws.Range("G15").Value = 0.002614821744523823
$ws.Range("K15").Value = 3.156011946184663
$ws.Range("L15").Value = 0.07969612643921664
$ws.Range("M15").Value = 0.8650991307610241

$ws.Range("C16").Value = 0.3580166647247154
$ws.Range("D16").Value = 0.4163131810062453
$ws.Range("E16").Value = 0.1177735675426916
$ws.Range("F16").Value = 5.976891006011812
$ws.Range("G16").Value = 0.002619426927467861
$ws.Range("K16").Value = 3.082197492872751
$ws.Range("L16").Value = 0.07942387701652365
$ws.Range("M16").Value = 0.8495960270680456

$ws.Range("C17").Value = 0.3558991712589261
$ws.Range("D17").Value = 0.4117709513040495
$ws.Range("E17").Value = 0.1173210596585896
$ws.Range("F17").Value = 5.885890223463861
$ws.Range("G17").Value = 0.002622310521792045
$ws.Range("K17").Value = 3.037660766039721
$ws.Range("L17").Value = 0.07926199354563934
$ws.Range("M17").Value = 0.8402652891950453

$ws.Range("C18").Value = 0.3546984010948506
$ws.Range("D18").Value = 0.4091621069620714
$ws.Range("E18").Value = 0.1170652186864878
$ws.Range("F18").Value = 5.83372683628852
$ws.Range("G18").Value = 0.002623990646744329
$ws.Range("K18").Value = 3.012317100434984
$ws.Range("L18").Value = 0.07917076087854369
$ws.Range("M18").Value = 0.8349643043825097

$ws.Range("C19").Value = 0.3542947807866881
$ws.Range("D19").Value = 0.4082794169670763
$ws.Range("E19").Value = 0.1169793545092261
$ws.Range("F19").Value = 5.816095440460145
$ws.Range("G19").Value = 0.002624563217479049
$ws.Range("K19").Value = 3.003782864582831
$ws.Range("L19").Value = 0.07914019294924302
$ws.Range("M19").Value = 0.8331807557939044

$ws.Range("C20").Value = 0.3561228048594103
$ws.Range("D20").Value = 0.4122540889968604
$ws.Range("E20").Value = 0.1173687710659266
$ws.Range("F20").Value = 5.895558929407997
$ws.Range("G20").Value = 0.002622001328549008
$ws.Range("K20").Value = 3.042373521862601
$ws.Range("L20").Value = 0.07927903173236572
$ws.Range("M20").Value = 0.8412517443064615

$ws.Range("C21").Value = 0.3624701979298379
$ws.Range("D21").Value = 0.4256544548566126
$ws.Range("E21").Value = 0.1187302017936744
$ws.Range("F21").Value = 6.164693923261268
$ws.Range("G21").Value = 0.002613650792175905
$ws.Range("K21").Value = 3.175301266459087
$ws.Range("L21").Value = 0.07976800917849403
$ws.Range("M21").Value = 0.8691576051219556

$ws.Range("C22").Value = 0.366772901673329
$ws.Range("D22").Value = 0.4344500893988368
$ws.Range("E22").Value = 0.1196597883044639
$ws.Range("F22").Value = 6.342244691628537
$ws.Range("G22").Value = 0.002608383662142361
$ws.Range("K22").Value = 3.264637451409783
$ws.Range("L22").Value = 0.08010447932791465
$ws.Range("M22").Value = 0.8879883649318145

$ws.Range("C23").Value = 0.3644622717517052
$ws.Range("D23").Value = 0.4297520206940533
$ws.Range("E23").Value = 0.1191599894869668
$ws.Range("F23").Value = 6.247327519515011
$ws.Range("G23").Value = 0.002611177478971366
$ws.Range("K23").Value = 3.216730355553977
$ws.Range("L23").Value = 0.07992334471813933
$ws.Range("M23").Value = 0.8778834764183188

$ws.Range("C24").Value = 0.3560216484043508
$ws.Range("D24").Value = 0.4120356545723212
$ws.Range("E24").Value = 0.1173471873123084
$ws.Range("F24").Value = 5.891187228019987
$ws.Range("G24").Value = 0.002622141045356515
$ws.Range("K24").Value = 3.040242072062256
$ws.Range("L24").Value = 0.07927132305018603
$ws.Range("M24").Value = 0.8408055706066051

$ws.Range("C25").Value = 0.3475720834848772
$ws.Range("D25").Value = 0.3930654491684891
$ws.Range("E25").Value = 0.1155609248740355
$ws.Range("F25").Value = 5.513789613939281
$ws.Range("G25").Value = 0.002634792163977298
$ws.Range("K25").Value = 2.860277640194624
$ws.Range("L25").Value = 0.0786397299044026
$ws.Range("M25").Value = 0.8033232237378769

